$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 716-717, shifting existing rows 716-769 down to 718-771
$ws.Range("A716:A717").EntireRow.Insert()

# Populate new row 716 (Murcott / Primera, new reporting week)
$ws.Cells.Item(716, 1).Value = 8
$ws.Cells.Item(716, 2).Value = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(716, 3).Value = 'Coquimbo'
$ws.Cells.Item(716, 4).Value = 45021
$ws.Cells.Item(716, 5).Value = 4
$ws.Cells.Item(716, 6).Value = 'Fruta'
$ws.Cells.Item(716, 7).Value = 100102
$ws.Cells.Item(716, 8).Value = 'Cítricos'
$ws.Cells.Item(716, 9).Value = 100102004
$ws.Cells.Item(716, 10).Value = 'Mandarina'
$ws.Cells.Item(716, 11).Value = 'Murcott'
$ws.Cells.Item(716, 12).Value = 'Primera'
$ws.Cells.Item(716, 13).Value = 20
$ws.Cells.Item(716, 14).Value = 250000
$ws.Cells.Item(716, 15).Value = 260000
$ws.Cells.Item(716, 16).Value = 255000
$ws.Cells.Item(716, 17).Value = '$/bins (450 kilos)'
$ws.Cells.Item(716, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(716, 19).Value = 567
$ws.Cells.Item(716, 20).Value = 450

# Populate new row 717 (Murcott / Segunda, new reporting week)
$ws.Cells.Item(717, 1).Value = 8
$ws.Cells.Item(717, 2).Value = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(717, 3).Value = 'Coquimbo'
$ws.Cells.Item(717, 4).Value = 45021
$ws.Cells.Item(717, 5).Value = 4
$ws.Cells.Item(717, 6).Value = 'Fruta'
$ws.Cells.Item(717, 7).Value = 100102
$ws.Cells.Item(717, 8).Value = 'Cítricos'
$ws.Cells.Item(717, 9).Value = 100102004
$ws.Cells.Item(717, 10).Value = 'Mandarina'
$ws.Cells.Item(717, 11).Value = 'Murcott'
$ws.Cells.Item(717, 12).Value = 'Segunda'
$ws.Cells.Item(717, 13).Value = 16
$ws.Cells.Item(717, 14).Value = 210000
$ws.Cells.Item(717, 15).Value = 220000
$ws.Cells.Item(717, 16).Value = 215000
$ws.Cells.Item(717, 17).Value = '$/bins (450 kilos)'
$ws.Cells.Item(717, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(717, 19).Value = 478
$ws.Cells.Item(717, 20).Value = 450
